# Casos de prueba.xlsx - update "Caso2" worksheet:
#  - extend the "Buscar pedidos" test case text (row 12)
#  - replace the "Cargar un pedido" test case (row 13) with a new
#    "Eliminar superuser" test case
#  - clear out the now-unused trailing rows (14-18)
#  - move the selection down to the newly emptied rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Caso2")

# Row 12: append clarifying sentences to the action/result description.
$ws.Range("D12").Value = "El superuser podrá ingresar un fecha de entrega para ver solo los pedidos de dicha fecha. Al presionar buscar devolverá listado de los pedidos y podrá verlos o editarlos. Sin no hay pedidos, mostrara una mensaje indicándolo. Si no carga fecha y presiona buscar, tambien devuelve mensaje."
$ws.Range("E12").Value = "Busca correctamente los pedidos por fecha, en caso que no haya devuelve mensaje en pantalla. Busco sin cargar fecha y tambien aparece pantalla con mensaje aclaratorio."
$ws.Rows.Item(12).RowHeight = 135

# Row 13: replace with a new test case about eliminating the superuser.
$ws.Range("C13").Value = "Eliminar superuser"
$ws.Range("D13").Value = "Estará bloqueda la eliminación del superuser desde el sitio. Solamente desde el admin será posible."
$ws.Range("E13").Value = "Al ingresar a editar el perfil no se muestra la opción para eliminarlo."
$ws.Rows.Item(13).RowHeight = 45

# Rows 14-18 no longer hold test cases; clear their contents and
# let the row height fall back to the sheet default (AutoFit drops the
# explicit height instead of pinning a custom one).
$ws.Range("A14:F18").ClearContents()
$ws.Range("A14:F18").Rows.AutoFit()

# Update the view: scroll back to the top and select the newly blank rows.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A14:F18").Select()
